$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '68.082.48'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +1.44%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.264.89'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '586.30'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '184.10'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +3.74%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -1.18%  '
$ws.Range('E9').Value = '  +3.33%  '
$ws.Range('E10').Value = '  -0.59%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.416'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.12%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '3.836.85'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('E13').Value = '  +0.30%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '28.58'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.37%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '68.144.74'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.53%  '
$ws.Range('E16').Value = '  +2.39%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.270.54'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.21%  '
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.62'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.97%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '382.65'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +2.55%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.70'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.65%  '
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '71.34'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.44%  '
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.0000120'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.69%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.82'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.56%  '
$ws.Range('E27').Value = '  +2.98%  '
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('E29').Value = '  +0.37%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.73'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.21%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.24'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +5.72%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '22.93'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.95%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.998'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.27'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.11%  '
$ws.Range('E35').Value = '  +2.47%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '162.56'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.74%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.87'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.25%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.836'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.62%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.78'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +5.02%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '26.66'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.71%  '
$ws.Range('E41').Value = '  +5.12%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.60'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.65%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '41.44'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.16%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '347.90'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.80%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '25.45'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.51%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0688'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.17%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.646.59'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.75%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0284'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.16%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '32.15'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +5.08%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.103'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.77%  '
$ws.Range('B51').Value = 'ONDO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.999'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.17%  '
